# EIA Table 2.7.B monthly update: roll the reporting period forward from
# "October 2016" to "November 2016" by adding the November data row and
# refreshing the Year-to-Date / Rolling-12-Months summary rows beneath it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the subtitle to reference November instead of October.
$ws.Range("A2").Value = "by Sector, 2006-November 2016 (Billion Btus)"

# 2) Insert a new row above the old "Year to Date" header (row 53) to hold
#    the new November monthly figures. This pushes every row below it
#    (headers, yearly summaries, and the notes row) down by one, which
#    also shifts the existing merged header ranges automatically.
$ws.Rows.Item(53).Insert()

# 3) Give the new row the same look as the other monthly data rows
#    (copy number formats / borders / fonts from the October row above it).
$ws.Range("A52:F52").Copy()
$ws.Range("A53:F53").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 4) Populate the new November monthly row.
$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 78622
$ws.Range("C53").Value = 789
$ws.Range("D53").Value = 1415
$ws.Range("E53").Value = 321
$ws.Range("F53").Value = 76098

# 5) The "Year to Date" header is now on row 54; refresh the annual totals
#    beneath it (rows 55-57) to the updated Year-to-Date figures.
$ws.Range("B55").Value = 863475
$ws.Range("C55").Value = 7918
$ws.Range("D55").Value = 20327
$ws.Range("E55").Value = 3457
$ws.Range("F55").Value = 831772

$ws.Range("B56").Value = 863003
$ws.Range("C56").Value = 8432
$ws.Range("D56").Value = 17381
$ws.Range("E56").Value = 3373
$ws.Range("F56").Value = 833817

$ws.Range("B57").Value = 841695
$ws.Range("C57").Value = 7508
$ws.Range("D57").Value = 15474
$ws.Range("E57").Value = 3944
$ws.Range("F57").Value = 814770

# 6) The "Rolling 12 Months" header is now on row 58; update its label to
#    say "November" and refresh the figures in rows 59-60.
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

$ws.Range("B59").Value = 945872
$ws.Range("C59").Value = 9348
$ws.Range("D59").Value = 19316
$ws.Range("E59").Value = 3682
$ws.Range("F59").Value = 913526

$ws.Range("B60").Value = 922654
$ws.Range("C60").Value = 8426
$ws.Range("D60").Value = 17293
$ws.Range("E60").Value = 4286
$ws.Range("F60").Value = 892649
